$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Events")

$ws.Cells.Item(2, 1).Value = 'Get Outta My Way!'
$ws.Cells.Item(2, 2).Value = 1
$ws.Cells.Item(2, 3).Value = 'Raise ⚠ for each hex tile with at least two Characters. %n Raise ⚠ for each Entrance tile with at least two Characters.'
$ws.Cells.Item(2, 4).Value = 'figure_alerts_neighboring.svg'

$ws.Cells.Item(3, 1).Value = 'Get Outta My Way!'
$ws.Cells.Item(3, 2).Value = 2
$ws.Cells.Item(3, 3).Value = 'For each hex tile with at least 2 characters on it, increase the alert bar by 1. For Entrance tiles, 2 alerts.'
$ws.Cells.Item(3, 4).Value = 'figure_alerts_neighboring.svg'

$ws.Cells.Item(4, 1).Value = 'Caught in a Reflection'
$ws.Cells.Item(4, 2).Value = 3
$ws.Cells.Item(4, 3).Value = 'Raise ⚠ for each Live Camera adjacent to a character.'
$ws.Cells.Item(4, 4).Value = 'figure_alerts_adj_camera.svg'

$ws.Cells.Item(5, 1).Value = 'Suspicious Guards'
$ws.Cells.Item(5, 2).Value = 4
$ws.Cells.Item(5, 3).Value = 'Raise ⚠ for each Character adjacent to an Unsubuded Guard '
$ws.Cells.Item(5, 4).Value = 'figure_alerts_adj_guard.svg'

$ws.Cells.Item(6, 1).Value = 'Tripped Over Stuff'
$ws.Cells.Item(6, 2).Value = 5
$ws.Cells.Item(6, 3).Value = 'Raise ⚠ for each Character on a Subdued Guard.%nRaise ⚠ for each Character on a Disabled Camera.'
$ws.Cells.Item(6, 4).Value = 'figure_alerts_on_disabled.svg'

$ws.Cells.Item(7, 1).Value = 'System Reboot'
$ws.Cells.Item(7, 2).Value = 6
$ws.Cells.Item(7, 3).Value = 'Make every Disabled Camera Live. Raise ⚠ every time you do this.'
$ws.Cells.Item(7, 4).Value = 'figure_alerts_power_on.svg'

$ws.Cells.Item(8, 1).Value = 'Shut the Gate!'
$ws.Cells.Item(8, 2).Value = 7
$ws.Cells.Item(8, 3).Value = 'Close the next Gate, starting with "A". Remove the tile. If a Character is on it, they are Busted. %nRaise ⚠⚠.'
$ws.Cells.Item(8, 4).Value = 'figure_alerts_lockdown_two_alerts.svg'

$ws.Cells.Item(9, 1).Value = 'Increased Suspician'
$ws.Cells.Item(9, 2).Value = 8
$ws.Cells.Item(9, 3).Value = 'Raise ⚠ for each Character adjacent to an Unsubdued Guard. Also, Raise ⚠ for each Character adjacent to a Live Camera.'
$ws.Cells.Item(9, 4).Value = 'figure_alerts_adj_camera_guard.svg'

$ws.Cells.Item(10, 1).Value = 'Tripped Over Stuff'
$ws.Cells.Item(10, 2).Value = 9
$ws.Cells.Item(10, 3).Value = 'Raise ⚠ for each Character on a Subdued Guard.%nRaise ⚠ for each Character on a Disabled Camera.'
$ws.Cells.Item(10, 4).Value = 'figure_alerts_on_disabled.svg'

$ws.Cells.Item(11, 1).Value = 'System Reboot'
$ws.Cells.Item(11, 2).Value = 10
$ws.Cells.Item(11, 3).Value = 'Make every Disabled Camera Live. Raise ⚠ every time you do this.'
$ws.Cells.Item(11, 4).Value = 'figure_alerts_power_on.svg'

$ws.Cells.Item(12, 1).Value = 'I Thought I Had Him'
$ws.Cells.Item(12, 2).Value = 11
$ws.Cells.Item(12, 3).Value = 'Raise ⚠⚠.%nUn-Subdue every Subdued Guard that shares a tile with a Character.'
$ws.Cells.Item(12, 4).Value = 'figure_alerts_unsubdue_alerts.svg'

$ws.Cells.Item(13, 1).Value = 'Increased Suspician'
$ws.Cells.Item(13, 2).Value = 12
$ws.Cells.Item(13, 3).Value = 'Raise ⚠ for each Character adjacent to an Unsubdued Guard. Also, Raise ⚠ for each Character adjacent to a Live Camera.'
$ws.Cells.Item(13, 4).Value = 'figure_alerts_adj_camera_guard.svg'

$ws.Cells.Item(14, 1).Value = 'Shut the Gate!'
$ws.Cells.Item(14, 2).Value = 13
$ws.Cells.Item(14, 3).Value = 'Close the next Gate, starting with "A". Remove the tile. If a Character is on it, they are Busted. %nRaise ⚠⚠.'
$ws.Cells.Item(14, 4).Value = 'figure_alerts_lockdown_two_alerts.svg'

$ws.Cells.Item(15, 1).Value = 'Tripped Over Stuff'
$ws.Cells.Item(15, 2).Value = 14
$ws.Cells.Item(15, 3).Value = 'Raise ⚠ for each Character on a Subdued Guard.%nRaise ⚠ for each Character on a Disabled Camera.'
$ws.Cells.Item(15, 4).Value = 'figure_alerts_on_disabled.svg'

$ws.Cells.Item(16, 1).Value = 'Hey You!'
$ws.Cells.Item(16, 2).Value = 15
$ws.Cells.Item(16, 3).Value = 'Add a Guard from the supply to each tile that has a Character adjacent to at least one Live Camera.'
$ws.Cells.Item(16, 4).Value = 'figure_alerts_hey_you.svg'

$ws.Cells.Item(17, 1).Value = 'Increased Suspician'
$ws.Cells.Item(17, 2).Value = 16
$ws.Cells.Item(17, 3).Value = 'Raise ⚠ for each Character adjacent to an Unsubdued Guard. Also, Raise ⚠ for each Character adjacent to a Live Camera.'
$ws.Cells.Item(17, 4).Value = 'figure_alerts_adj_camera_guard.svg'

$ws.Cells.Item(18, 1).Value = 'Shut the Gate!'
$ws.Cells.Item(18, 2).Value = 17
$ws.Cells.Item(18, 3).Value = 'Close the next Gate, starting with "A". Remove the tile. If a Character is on it, they are Busted. %nRaise ⚠.'
$ws.Cells.Item(18, 4).Value = 'figure_alerts_lockdown_one_alert.svg'

$ws.Cells.Item(19, 1).Value = 'Ties Broke!'
$ws.Cells.Item(19, 2).Value = 18
$ws.Cells.Item(19, 3).Value = 'Un-Subdue every Subdued Guard that shares a tile with a Character.'
$ws.Cells.Item(19, 4).Value = 'figure_alerts_unsubdue.svg'

$ws.Cells.Item(20, 1).Value = 'Go Check on It'
$ws.Cells.Item(20, 2).Value = 19
$ws.Cells.Item(20, 3).Value = 'Replace every Disabled Camera with a Guard.'
$ws.Cells.Item(20, 4).Value = 'figure_alerts_call_it_in.svg'

$ws.Cells.Item(21, 1).Value = 'Shut the Gates!'
$ws.Cells.Item(21, 2).Value = 20
$ws.Cells.Item(21, 3).Value = 'Close the next Gate, starting with "A". Remove the tile. If a character is on it, they are immediately Busted. '
$ws.Cells.Item(21, 4).Value = 'figure_alerts_lockdown.svg'

$ws.Cells.Item(22, 1).Value = 'Sound the Alarm!'
$ws.Cells.Item(22, 2).Value = 21
$ws.Cells.Item(22, 3).Value = 'Initiate Escape Phase.'
$ws.Cells.Item(22, 4).Value = 'figure_alerts_alarm.svg'

$ws.Activate() | Out-Null
$ws.Range("C16").Select() | Out-Null